# Update the "aca-metal-tier" StructureDefinition workbook:
#  - bump the canonical URL from ibm.com to linuxforhealth.org
#  - bump Version 7.0.0 -> 8.0.0
#  - bump the publication Date
#  - rename the Publisher from "Alvearie Team" to "LinuxForHealth Team"
#  - clear the stale Constraint(s) text that had been left on the
#    Extension row (it now only belongs on Extension.extension)

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/aca-metal-tier"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/aca-metal-tier"
